$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.9029038349787394
$ws.Cells.Item(2, 4).Value = 0.9167482256889343
$ws.Cells.Item(2, 5).Value = 0.8286526948213577
$ws.Cells.Item(2, 6).Value = 0.8402521014213562
$ws.Cells.Item(2, 7).Value = 0.8015593141317368
$ws.Cells.Item(2, 8).Value = 0.9777853488922119
$ws.Cells.Item(2, 9).Value = 0.9856123328208923
$ws.Cells.Item(2, 10).Value = 0.7952601909637451
$ws.Cells.Item(2, 11).Value = 0.9316011667251587
$ws.Cells.Item(3, 2).Value = 0.8623180290063223
$ws.Cells.Item(3, 4).Value = 0.756869912147522
$ws.Cells.Item(3, 5).Value = 0.830847313006719
$ws.Cells.Item(3, 6).Value = 0.8537154396375021
$ws.Cells.Item(3, 7).Value = 0.8012186884880066
$ws.Cells.Item(3, 8).Value = 0.9441962639490763
$ws.Cells.Item(3, 9).Value = 0.9852280616760254
$ws.Cells.Item(3, 10).Value = 0.7365642786026001
$ws.Cells.Item(3, 11).Value = 0.9521704117457072
$ws.Cells.Item(4, 2).Value = 0.9504640698432922
$ws.Cells.Item(4, 3).Value = 0.8124158382415771
$ws.Cells.Item(4, 5).Value = 0.9371106326580048
$ws.Cells.Item(4, 6).Value = 0.9950271844863892
$ws.Cells.Item(4, 7).Value = 0.9532506465911865
$ws.Cells.Item(4, 8).Value = 0.95600825548172
$ws.Cells.Item(4, 9).Value = 0.9846462905406952
$ws.Cells.Item(4, 10).Value = 0.9648858904838562
$ws.Cells.Item(4, 11).Value = 0.9573302268981934
$ws.Cells.Item(5, 2).Value = 0.8235698044300079
$ws.Cells.Item(5, 3).Value = 0.821491003036499
$ws.Cells.Item(5, 4).Value = 0.9481375217437744
$ws.Cells.Item(5, 6).Value = 0.9241468608379364
$ws.Cells.Item(5, 7).Value = 0.8520738482475281
$ws.Cells.Item(5, 8).Value = 0.9246824085712433
$ws.Cells.Item(5, 9).Value = 0.9620617628097534
$ws.Cells.Item(5, 10).Value = 0.9240339994430542
$ws.Cells.Item(5, 11).Value = 0.8459548950195312
$ws.Cells.Item(6, 2).Value = 0.8665298223495483
$ws.Cells.Item(6, 3).Value = 0.7938462098439535
$ws.Cells.Item(6, 4).Value = 0.9981905817985535
$ws.Cells.Item(6, 5).Value = 0.9416531324386597
$ws.Cells.Item(6, 7).Value = 0.9633549749851227
$ws.Cells.Item(6, 8).Value = 0.9462270736694336
$ws.Cells.Item(6, 9).Value = 0.9797036647796631
$ws.Cells.Item(6, 10).Value = 0.9467079043388367
$ws.Cells.Item(6, 11).Value = 0.9381299614906311
$ws.Cells.Item(7, 2).Value = 0.8539272546768188
$ws.Cells.Item(7, 3).Value = 0.8183208306630453
$ws.Cells.Item(7, 4).Value = 0.9335795938968658
$ws.Cells.Item(7, 5).Value = 0.8379693478345871
$ws.Cells.Item(7, 6).Value = 0.9150501489639282
$ws.Cells.Item(7, 8).Value = 0.9320261478424072
$ws.Cells.Item(7, 9).Value = 0.9759467393159866
$ws.Cells.Item(7, 10).Value = 0.865458756685257
$ws.Cells.Item(7, 11).Value = 0.8839512765407562
$ws.Cells.Item(8, 2).Value = 0.973026305437088
$ws.Cells.Item(8, 3).Value = 0.9541288812955221
$ws.Cells.Item(8, 4).Value = 0.9361419677734375
$ws.Cells.Item(8, 5).Value = 0.9305903613567352
$ws.Cells.Item(8, 6).Value = 0.9607416987419128
$ws.Cells.Item(8, 7).Value = 0.9218071103096008
$ws.Cells.Item(8, 9).Value = 0.9983170628547668
$ws.Cells.Item(8, 10).Value = 0.7977302670478821
$ws.Cells.Item(8, 11).Value = 0.9633198976516724
$ws.Cells.Item(9, 2).Value = 0.9920817166566849
$ws.Cells.Item(9, 3).Value = 0.9894202649593353
$ws.Cells.Item(9, 4).Value = 0.9839665591716766
$ws.Cells.Item(9, 5).Value = 0.9576843231916428
$ws.Cells.Item(9, 6).Value = 0.9884465038776398
$ws.Cells.Item(9, 7).Value = 0.9761862307786942
$ws.Cells.Item(9, 8).Value = 0.9986096322536469
$ws.Cells.Item(9, 10).Value = 0.9546948075294495
$ws.Cells.Item(9, 11).Value = 0.9956967830657959
$ws.Cells.Item(10, 2).Value = 0.8155481219291687
$ws.Cells.Item(10, 3).Value = 0.6833437085151672
$ws.Cells.Item(10, 4).Value = 0.9411723613739014
$ws.Cells.Item(10, 5).Value = 0.8764290809631348
$ws.Cells.Item(10, 6).Value = 0.9361354112625122
$ws.Cells.Item(10, 7).Value = 0.902280867099762
$ws.Cells.Item(10, 8).Value = 0.8584713935852051
$ws.Cells.Item(10, 9).Value = 0.9764431118965149
$ws.Cells.Item(10, 11).Value = 0.8363800644874573
$ws.Cells.Item(11, 2).Value = 0.8516279458999634
$ws.Cells.Item(11, 3).Value = 0.9452718098958334
$ws.Cells.Item(11, 4).Value = 0.9540036916732788
$ws.Cells.Item(11, 5).Value = 0.7847095429897308
$ws.Cells.Item(11, 6).Value = 0.9205443859100342
$ws.Cells.Item(11, 7).Value = 0.8490531444549561
$ws.Cells.Item(11, 8).Value = 0.95635986328125
$ws.Cells.Item(11, 9).Value = 0.9979508817195892
$ws.Cells.Item(11, 10).Value = 0.766191840171814
